# New datasets in 3 different formats
# Updates the measured/computed values in columns B (u$t (s)$), C ($V (V)$)
# and D (u$V (V)$) for rows 2-26 on Sheet1 with new dataset values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.02548244374480401
$ws.Range("C2").Value = 155.1196868684945
$ws.Range("D2").Value = 4.425198365858968

$ws.Range("B3").Value = 0.05535919592979506
$ws.Range("C3").Value = 206.4080325671268
$ws.Range("D3").Value = 1.731032935600061

$ws.Range("B4").Value = 0.02522128400978997
$ws.Range("C4").Value = 210.5360660892276
$ws.Range("D4").Value = 3.736248590126293

$ws.Range("B5").Value = 0.04132126513780501
$ws.Range("C5").Value = 155.3305724094011
$ws.Range("D5").Value = 3.7838069549992

$ws.Range("B6").Value = 0.03466915599985192
$ws.Range("C6").Value = 53.17785819930636
$ws.Range("D6").Value = 4.160981177509423

$ws.Range("B7").Value = 0.05248196264365405
$ws.Range("C7").Value = -56.80209460711851
$ws.Range("D7").Value = 3.311827795540527

$ws.Range("B8").Value = 0.06199905553525949
$ws.Range("C8").Value = -149.7835605811072
$ws.Range("D8").Value = 2.578273756700676

$ws.Range("B9").Value = 0.0307042228034687
$ws.Range("C9").Value = -208.4610238095932
$ws.Range("D9").Value = 2.60837448200821

$ws.Range("B10").Value = 0.04351143443474575
$ws.Range("C10").Value = -208.1500919880107
$ws.Range("D10").Value = 3.439300039355125

$ws.Range("B11").Value = 0.06660062446966006
$ws.Range("C11").Value = -158.9717853950779
$ws.Range("D11").Value = 4.205167506484173

$ws.Range("B12").Value = 0.05262879062982387
$ws.Range("C12").Value = -56.88293046092694
$ws.Range("D12").Value = 1.886673272285295

$ws.Range("B13").Value = 0.04903752791106657
$ws.Range("C13").Value = 66.57224432659994
$ws.Range("D13").Value = 2.03434059635256

$ws.Range("B14").Value = 0.05631588178261286
$ws.Range("C14").Value = 146.5790142715803
$ws.Range("D14").Value = 2.18753896483895

$ws.Range("B15").Value = 0.05480028572849174
$ws.Range("C15").Value = 212.9642546300508
$ws.Range("D15").Value = 2.556062799602127

$ws.Range("B16").Value = 0.0681480504269887
$ws.Range("C16").Value = 215.3856117547747
$ws.Range("D16").Value = 1.878672372521283

$ws.Range("B17").Value = 0.04809261621359028
$ws.Range("C17").Value = 150.449981511342
$ws.Range("D17").Value = 4.345239357847176

$ws.Range("B18").Value = 0.06857243578149287
$ws.Range("C18").Value = 57.21635650895196
$ws.Range("D18").Value = 1.524070364759401

$ws.Range("B19").Value = 0.0618898777405417
$ws.Range("C19").Value = -61.13079152206469
$ws.Range("D19").Value = 1.909271003074378

$ws.Range("B20").Value = 0.0461159512399767
$ws.Range("C20").Value = -151.2753811119431
$ws.Range("D20").Value = 1.675657304973546

$ws.Range("B21").Value = 0.04750187297072375
$ws.Range("C21").Value = -215.9964746451515
$ws.Range("D21").Value = 3.579205372380471

$ws.Range("B22").Value = 0.03677497440793234
$ws.Range("C22").Value = -207.8761722410015
$ws.Range("D22").Value = 4.038038653528557

$ws.Range("B23").Value = 0.04691562057841696
$ws.Range("C23").Value = -148.3765779643703
$ws.Range("D23").Value = 3.476279691001706

$ws.Range("B24").Value = 0.05920284380409101
$ws.Range("C24").Value = -68.24650357824092
$ws.Range("D24").Value = 3.434737148204446

$ws.Range("B25").Value = 0.05275404804077709
$ws.Range("C25").Value = 60.53144143265828
$ws.Range("D25").Value = 2.32032468974526

$ws.Range("B26").Value = 0.03641645563398804
$ws.Range("C26").Value = 151.8903945183057
$ws.Range("D26").Value = 1.984106985647873
